$wb = $excel.ActiveWorkbook

# Rename "Sheet1" to "RestAssured" and populate it with the Rest Assured
# data-driven-test fixture data.
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "RestAssured"

# Header row
$ws.Range("A1").Value = "Course"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "isbn"
$ws.Range("D1").Value = "aisle"
$ws.Range("E1").Value = "name"

# Data row
$ws.Range("A2").Value = "RestAPI"
$ws.Range("B2").Value = "Learn Appium Automation with Java"
$ws.Range("C2").Value = "fhhv"
$ws.Range("D2").Value = 637628
$ws.Range("E2").Value = "Suriya Kumarr"

# Style the data row (except the numeric aisle cell) with a small teal
# Consolas font. Format A2 directly, then fan the same look out to the
# other text cells in the row via a format-only paste so they all share
# one cell style.
$ws.Range("A2").Font.Name = "Consolas"
$ws.Range("A2").Font.Family = 3
$ws.Range("A2").Font.Size = 10
$ws.Range("A2").Font.Color = 10733079

$ws.Range("A2").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("E2").PasteSpecial(-4122)

# Leave the TestData sheet's header row selected (but not the active
# sheet/tab).
$tdws = $wb.Worksheets.Item("TestData")
[void]$tdws.Range("A1:D1").Select()

# Select C2 on the new sheet and make it the active tab.
[void]$ws.Range("C2").Select()
$ws.Activate()
